$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 71; this shifts the existing rows 71-82 down to 72-83.
$ws.Rows.Item(71).Insert()

# Populate the newly inserted row 71 with the new data record.
$ws.Range("A71").Value = 11
$ws.Range("B71").Value = "Vega Monumental Concepción"
$ws.Range("C71").Value = "Bíobío"
$ws.Range("D71").Value = 44505
$ws.Range("E71").Value = 8
$ws.Range("F71").Value = 100112043
$ws.Range("G71").Value = "Pepino ensalada"
$ws.Range("H71").Value = "Sin especificar"
$ws.Range("I71").Value = "Primera"
$ws.Range("J71").Value = 300
$ws.Range("K71").Value = 6500
$ws.Range("L71").Value = 7000
$ws.Range("M71").Value = 6750
$ws.Range("N71").Value = "$/caja 80 unidades"
$ws.Range("O71").Value = "Región del Maule"
$ws.Range("P71").Value = 84
$ws.Range("Q71").Value = 80
$ws.Range("R71").Value = "Hortaliza"

# Match the date-format style used by the other rows in column D.
$ws.Range("D71").NumberFormat = "YYYY-MM-DD HH:MM:SS"
